$wb = $excel.ActiveWorkbook

# 1) Rename the month sheets to "M_YYYY" form.
$renames = @{
    "October"   = "10_2016"
    "September" = "9_2016"
    "August"    = "8_2016"
    "July"      = "7_2016"
    "June"      = "6_2016"
    "May"       = "5_2016"
    "April"     = "4_2016"
    "March"     = "3_2016"
    "February"  = "2_2016"
    "January"   = "1_2016"
    "December"  = "12_2015"
    "November"  = "11_2015"
}

foreach ($ws in $wb.Worksheets) {
    $newName = $renames[$ws.Name]
    if ($newName) {
        $ws.Name = $newName
    }
}

# 2) Re-write the header row on every sheet with the new column labels
#    (GMLID stays, the 4 unit columns get renamed from "X (unit)" to "X_unit").
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Value = "GMLID"
    $ws.Range("B1").Value = "Electricity_KWH"
    $ws.Range("C1").Value = "Cold_Water_m3"
    $ws.Range("D1").Value = "Hot_Water_m3"
    $ws.Range("E1").Value = "Heat_MWH"
}

# 3) Widen column B on the "6_2016" (June) sheet to fit its header text.
$wb.Worksheets.Item("6_2016").Columns.Item(2).ColumnWidth = 17.140625

# 4) Give "1_2016" (January) and "12_2015" (December) sheets a page setup
#    matching the other sheets (paper size 9 / portrait).
foreach ($name in @("1_2016", "12_2015", "11_2015")) {
    $ps = $wb.Worksheets.Item($name).PageSetup
    $ps.PaperSize = 9
    $ps.Orientation = 1
}

# 5) Reset every sheet's selection to the header row, then restore the
#    previously-active / newly-active tab state.
foreach ($ws in $wb.Worksheets) {
    $ws.Activate() | Out-Null
    $ws.Range("A1:E1").Select() | Out-Null
}

$wb.Worksheets.Item("11_2015").Range("E1").Select() | Out-Null

$wb.Worksheets.Item("10_2016").Activate() | Out-Null
